# Slide 7 ("When to not use Redis") - expand on the SQL Server queries bullet
# and clarify why Redis is still useful alongside SQL Server.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Locate the body content placeholder by name (avoids relying on a fixed
# shape index).
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -like "Content Placeholder*") {
        $shape = $candidate
        break
    }
}

$tf = $shape.TextFrame
$tr = $tf.TextRange

# Find the "SQL Server Queries ..." bullet and the trailing "Redis performs
# at about the same speed" bullet by matching their text, rather than
# assuming fixed paragraph indices. Paragraph.Text carries a trailing
# carriage return (except for the very last paragraph), so trim before
# comparing.
$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $paraText = $para.Text.TrimEnd("`r")
    if ($paraText -eq "SQL Server Queries (Specifically the results)") {
        $run = $para.Runs(1, 1)
        $run.Text = "SQL Server Queries (Specifically the results) for performance improvements"
    }
}

$lastPara = $tr.Paragraphs($paraCount, 1)
$lastParaText = $lastPara.Text.TrimEnd("`r")
if ($lastParaText -eq "Redis performs at about the same speed") {
    # Append a new sub-bullet (same indent level as the paragraph it
    # follows) after the final paragraph in the placeholder.
    $null = $tr.InsertAfter("`rOnly reason would be to remove load from SQL Server")
}
